$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 159.71428
$ws.Range("I2").Value = 64.75
$ws.Range("K2").Value = 64.75
$ws.Range("M2").Value = 48.25

$ws.Range("H28").Value = 94544.17999999999
$ws.Range("I28").Value = 126246
$ws.Range("K28").Value = 126246
$ws.Range("M28").Value = -125761

$ws.Range("H33").Value = 195.41667
$ws.Range("I33").Value = 199.63637
$ws.Range("K33").Value = 199.63637
$ws.Range("M33").Value = 29.36363

$ws.Range("H69").Value = 7513.3335
$ws.Range("I69").Value = 5000
$ws.Range("J69").Value = 8770
$ws.Range("K69").Value = 15000
$ws.Range("L69").Value = 26310
$ws.Range("M69").Value = -14126
$ws.Range("N69").Value = -28058

$ws.Range("H72").Value = 7513.3335
$ws.Range("I72").Value = 5000
$ws.Range("J72").Value = 8770
$ws.Range("K72").Value = 45000
$ws.Range("L72").Value = 78930
$ws.Range("M72").Value = -40632
$ws.Range("N72").Value = -87666

$ws.Range("H86").Value = 6205.6665
$ws.Range("I86").Value = 5907.2856
$ws.Range("K86").Value = 5907.2856
$ws.Range("M86").Value = -4784.2856

$ws.Range("H88").Value = 12995.454
$ws.Range("J88").Value = 16619.5
$ws.Range("L88").Value = 16619.5
$ws.Range("N88").Value = -17431.5

$ws.Range("H89").Value = 6205.6665
$ws.Range("I89").Value = 5907.2856
$ws.Range("K89").Value = 29536.428
$ws.Range("M89").Value = -23920.428

$ws.Range("H91").Value = 12995.454
$ws.Range("J91").Value = 16619.5
$ws.Range("L91").Value = 16619.5
$ws.Range("N91").Value = -19427.5

$ws.Range("H98").Value = 2708.6667
$ws.Range("I98").Value = 502.92307
$ws.Range("K98").Value = 502.92307
$ws.Range("M98").Value = 995.0769299999999

$ws.Range("H101").Value = 148.5
$ws.Range("J101").Value = 0
$ws.Range("L101").Value = 0
$ws.Range("N101").ClearContents()

$ws.Range("H122").Value = 2708.6667
$ws.Range("I122").Value = 502.92307
$ws.Range("K122").Value = 1508.76921
$ws.Range("M122").Value = 941.2307900000001

$ws.Range("H132").Value = 1821.9048
$ws.Range("I132").Value = 1821.9048
$ws.Range("K132").Value = 5465.7144
$ws.Range("M132").Value = -2935.7144

$ws.Range("H137").Value = 4377.6665
$ws.Range("J137").Value = 4377.6665
$ws.Range("L137").Value = 13132.9995
$ws.Range("N137").Value = -18232.9995

$ws.Range("H138").Value = 3302.6287
$ws.Range("J138").Value = 3155.76
$ws.Range("L138").Value = 9467.280000000001
$ws.Range("N138").Value = -19747.28

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 3086.3057
$ws.Range("I61").Value = 2838.9644
$ws.Range("K61").Value = 2838.9644
$ws.Range("M61").Value = -2626.9644

$ws.Range("H122").Value = 2162.1667
$ws.Range("I122").Value = 1701.1177
$ws.Range("K122").Value = 5103.3531
$ws.Range("M122").Value = -2653.3531

$ws.Range("H132").Value = 5238.533
$ws.Range("I132").Value = 3119.2307
$ws.Range("J132").Value = 19014
$ws.Range("K132").Value = 9357.6921
$ws.Range("L132").Value = 57042
$ws.Range("M132").Value = -6827.6921
$ws.Range("N132").Value = -62102

$ws.Range("H136").Value = 3086.3057
$ws.Range("I136").Value = 2838.9644
$ws.Range("K136").Value = 8516.893199999999
$ws.Range("M136").Value = -5966.893199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 3212.1428
$ws.Range("I134").Value = 1705
$ws.Range("J134").Value = 6980
$ws.Range("K134").Value = 5115
$ws.Range("L134").Value = 20940
$ws.Range("M134").Value = -2580
$ws.Range("N134").Value = -26010

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 29201.56
$ws.Range("I31").Value = 3097.6072
$ws.Range("J31").Value = 85425.46000000001
$ws.Range("K31").Value = 3097.6072
$ws.Range("L31").Value = 85425.46000000001
$ws.Range("M31").Value = -2802.6072
$ws.Range("N31").Value = -86015.46000000001

$ws.Range("H34").Value = 29201.56
$ws.Range("I34").Value = 3097.6072
$ws.Range("J34").Value = 85425.46000000001
$ws.Range("K34").Value = 3097.6072
$ws.Range("L34").Value = 85425.46000000001
$ws.Range("M34").Value = -2895.6072
$ws.Range("N34").Value = -85829.46000000001

$ws.Range("H58").Value = 5884.8184
$ws.Range("I58").Value = 1463.625
$ws.Range("K58").Value = 1463.625
$ws.Range("M58").Value = -1260.625

$ws.Range("H132").Value = 4860.2407
$ws.Range("I132").Value = 4682.755
$ws.Range("K132").Value = 14048.265
$ws.Range("M132").Value = -11518.265

$ws.Range("H136").Value = 5884.8184
$ws.Range("I136").Value = 1463.625
$ws.Range("K136").Value = 4390.875
$ws.Range("M136").Value = -1840.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 2432
$ws.Range("I34").Value = 1759.4286
$ws.Range("J34").Value = 3216.6667
$ws.Range("K34").Value = 5278.2858
$ws.Range("L34").Value = 9650.000100000001
$ws.Range("M34").Value = -5194.2858
$ws.Range("N34").Value = -9818.000100000001

$ws.Range("H39").Value = 3309.4546
$ws.Range("I39").Value = 2134
$ws.Range("K39").Value = 6402
$ws.Range("M39").Value = -6108

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H33").Value = 18999
$ws.Range("J33").Value = 18999
$ws.Range("L33").Value = 18999
$ws.Range("N33").Value = -19503

$ws.Range("H97").Value = 1099
$ws.Range("I97").Value = 885.9
$ws.Range("K97").Value = 885.9
$ws.Range("M97").Value = -389.9

$ws.Range("H100").Value = 0
$ws.Range("J100").Value = 0
$ws.Range("L100").Value = 0
$ws.Range("N100").ClearContents()

$ws.Range("H102").Value = 2228.5715
$ws.Range("I102").Value = 1430.8948
$ws.Range("K102").Value = 1430.8948
$ws.Range("M102").Value = 191.1052

$ws.Range("H122").Value = 2710.6667
$ws.Range("I122").Value = 2217
$ws.Range("K122").Value = 6651
$ws.Range("M122").Value = -4201

$ws.Range("H126").Value = 3089.1765
$ws.Range("I126").Value = 2000.1538
$ws.Range("J126").Value = 6628.5
$ws.Range("K126").Value = 6000.4614
$ws.Range("L126").Value = 19885.5
$ws.Range("M126").Value = -3530.4614
$ws.Range("N126").Value = -24825.5

$ws.Range("H132").Value = 175335.5
$ws.Range("I132").Value = 335999.66
$ws.Range("J132").Value = 14671.333
$ws.Range("K132").Value = 1007998.98
$ws.Range("L132").Value = 44013.999
$ws.Range("M132").Value = -1005468.98
$ws.Range("N132").Value = -49073.999

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 4208.778
$ws.Range("I16").Value = 3672.375
$ws.Range("K16").Value = 3672.375
$ws.Range("M16").Value = -3502.375

$ws.Range("H22").Value = 4375.0835
$ws.Range("I22").Value = 2100
$ws.Range("J22").Value = 6000.143
$ws.Range("K22").Value = 2100
$ws.Range("L22").Value = 6000.143
$ws.Range("M22").Value = -1805
$ws.Range("N22").Value = -6590.143

$ws.Range("H27").Value = 4375.0835
$ws.Range("I27").Value = 2100
$ws.Range("J27").Value = 6000.143
$ws.Range("K27").Value = 2100
$ws.Range("L27").Value = 6000.143
$ws.Range("M27").Value = -1993
$ws.Range("N27").Value = -6214.143

$ws.Range("H46").Value = 4374.375
$ws.Range("I46").Value = 1498.5
$ws.Range("K46").Value = 1498.5
$ws.Range("M46").Value = -1310.5

$ws.Range("H82").Value = 4725.421
$ws.Range("J82").Value = 7943.6665
$ws.Range("L82").Value = 7943.6665
$ws.Range("N82").Value = -8665.666499999999

$ws.Range("H85").Value = 4725.421
$ws.Range("J85").Value = 7943.6665
$ws.Range("L85").Value = 7943.6665
$ws.Range("N85").Value = -10439.6665

$ws.Range("H93").Value = 2878.0557
$ws.Range("I93").Value = 2283
$ws.Range("J93").Value = 3621.875
$ws.Range("K93").Value = 2283
$ws.Range("L93").Value = 3621.875
$ws.Range("M93").Value = -1035
$ws.Range("N93").Value = -6117.875

$ws.Range("H100").Value = 5823.75
$ws.Range("I100").Value = 1430.3334
$ws.Range("K100").Value = 1430.3334
$ws.Range("M100").Value = -889.3334

$ws.Range("H132").Value = 4662.8184
$ws.Range("I132").Value = 2259
$ws.Range("J132").Value = 6666
$ws.Range("K132").Value = 6777
$ws.Range("L132").Value = 19998
$ws.Range("M132").Value = -4247
$ws.Range("N132").Value = -25058

$ws.Range("H136").Value = 5881.2964
$ws.Range("I136").Value = 2649.0667
$ws.Range("J136").Value = 9921.583000000001
$ws.Range("K136").Value = 7947.2001
$ws.Range("L136").Value = 29764.749
$ws.Range("M136").Value = -5397.2001
$ws.Range("N136").Value = -34864.749

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 4155.3335
$ws.Range("J81").Value = 19002
$ws.Range("L81").Value = 38004
$ws.Range("N81").Value = -40126

$ws.Range("H84").Value = 4155.3335
$ws.Range("J84").Value = 19002
$ws.Range("L84").Value = 190020
$ws.Range("N84").Value = -200628

$ws.Range("H136").Value = 11073
$ws.Range("I136").Value = 6816.3335
$ws.Range("J136").Value = 15329.667
$ws.Range("K136").Value = 20449.0005
$ws.Range("L136").Value = 45989.001
$ws.Range("M136").Value = -17899.0005
$ws.Range("N136").Value = -51089.001
